$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 333
$ws.Range("E2").Value = 69
$ws.Range("F2").Value = 69
$ws.Range("G2").Value = 31
$ws.Range("H2").Value = 28
$ws.Range("I2").Value = 31
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 1786
$ws.Range("L2").Value = 1034
$ws.Range("M2").Value = 752
$ws.Range("N2").Value = 800
$ws.Range("O2").Value = -48
$ws.Range("P2").Value = 351
$ws.Range("Q2").Value = 129
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = -98
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 126
$ws.Range("V2").Value = 464
$ws.Range("W2").Value = 20.72
$ws.Range("X2").Value = 8.29
$ws.Range("Y2").Value = 3.99
$ws.Range("Z2").Value = 1.58
$ws.Range("AA2").Value = 137.47
$ws.Range("AB2").Value = 119.5
$ws.Range("AC2").Value = 442
$ws.Range("AD2").Value = 8.6
$ws.Range("AE2").Value = 11393
$ws.Range("AF2").Value = 0.33
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 7022839

# Row 3
$ws.Range("D3").Value = 76
$ws.Range("E3").Value = -19
$ws.Range("F3").Value = -19
$ws.Range("G3").Value = -31
$ws.Range("H3").Value = -15
$ws.Range("I3").Value = -12
$ws.Range("J3").Value = -3
$ws.Range("K3").Value = 1735
$ws.Range("L3").Value = 1044
$ws.Range("M3").Value = 691
$ws.Range("N3").Value = 741
$ws.Range("O3").Value = -51
$ws.Range("P3").Value = 351
$ws.Range("Q3").Value = -23
$ws.Range("R3").Value = 48
$ws.Range("S3").Value = -58
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = -25
$ws.Range("V3").Value = 404
$ws.Range("W3").Value = -24.45
$ws.Range("X3").Value = -20.15
$ws.Range("Y3").Value = -1.54
$ws.Range("Z3").Value = -0.87
$ws.Range("AA3").Value = 151.19
$ws.Range("AB3").Value = 98.38
$ws.Range("AC3").Value = -169
$ws.Range("AD3").Value = -36.75
$ws.Range("AE3").Value = 10557
$ws.Range("AF3").Value = 0.59
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 7022839

# Row 4
$ws.Range("D4").Value = 58
$ws.Range("E4").Value = -91
$ws.Range("F4").Value = -91
$ws.Range("G4").Value = -125
$ws.Range("H4").Value = -172
$ws.Range("I4").Value = -165
$ws.Range("J4").Value = -7
$ws.Range("K4").Value = 1750
$ws.Range("L4").Value = 1229
$ws.Range("M4").Value = 521
$ws.Range("N4").Value = 578
$ws.Range("O4").Value = -58
$ws.Range("P4").Value = 351
$ws.Range("Q4").Value = -84
$ws.Range("R4").Value = -42
$ws.Range("S4").Value = 192
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = -85
$ws.Range("V4").Value = 597
$ws.Range("W4").Value = -156.56
$ws.Range("X4").Value = -294.25
$ws.Range("Y4").Value = -24.95
$ws.Range("Z4").Value = -9.86
$ws.Range("AA4").Value = 236.08
$ws.Range("AB4").Value = 51.51
$ws.Range("AC4").Value = -2344
$ws.Range("AD4").Value = -2.3
$ws.Range("AE4").Value = 8237
$ws.Range("AF4").Value = 0.65
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 7022839
$ws.Range("AG4:AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 506
$ws.Range("E5").Value = 97
$ws.Range("F5").Value = 97
$ws.Range("G5").Value = -39
$ws.Range("H5").Value = -29
$ws.Range("I5").Value = -27
$ws.Range("J5").Value = -2
$ws.Range("K5").Value = 1522
$ws.Range("L5").Value = 1028
$ws.Range("M5").Value = 494
$ws.Range("N5").Value = 553
$ws.Range("O5").Value = -60
$ws.Range("P5").Value = 351
$ws.Range("Q5").Value = 141
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = -139
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 141
$ws.Range("V5").Value = 458
$ws.Range("W5").Value = 19.19
$ws.Range("X5").Value = -5.65
$ws.Range("Y5").Value = -4.73
$ws.Range("Z5").Value = -1.75
$ws.Range("AA5").Value = 208.36
$ws.Range("AB5").Value = 43.89
$ws.Range("AC5").Value = -381
$ws.Range("AD5").Value = -12.94
$ws.Range("AE5").Value = 7877
$ws.Range("AF5").Value = 0.63
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 7022839
$ws.Range("AG5:AH5").ClearContents()

# Row 6
$ws.Range("D6").Value = 580
$ws.Range("E6").Value = 230
$ws.Range("F6").Value = 230
$ws.Range("G6").Value = 132
$ws.Range("H6").Value = -16
$ws.Range("I6").Value = -14
$ws.Range("K6").Value = 1617
$ws.Range("L6").Value = 1209
$ws.Range("M6").Value = 408
$ws.Range("N6").Value = 469
$ws.Range("P6").Value = 351
$ws.Range("Q6").Value = -1
$ws.Range("R6").Value = -13
$ws.Range("S6").Value = 14
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = -1
$ws.Range("V6").Value = 438
$ws.Range("W6").Value = 39.67
$ws.Range("X6").Value = -2.76
$ws.Range("Y6").Value = -2.77
$ws.Range("Z6").Value = -1.02
$ws.Range("AA6").Value = 296.33
$ws.Range("AB6").Value = 27.02
$ws.Range("AC6").Value = -202
$ws.Range("AD6").Value = -16.94
$ws.Range("AE6").Value = 6682
$ws.Range("AF6").Value = 0.51
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 7022839
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9: clear D:AJ (only A,B,C retained)
$ws.Range("D7:AJ9").ClearContents()
